$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '70.189.86'
$ws.Range('E2').Value = '  -0.64%  '

$ws.Range('D3').Value = '3.509.39'
$ws.Range('E3').Value = '  -0.90%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').Value = '''604.90'
$ws.Range('E5').Value = '  -0.54%  '

$ws.Range('D6').Value = '''172.89'
$ws.Range('E6').Value = '  -0.84%  '

$ws.Range('D7').Value = '''0.608'
$ws.Range('E7').Value = '  -1.61%  '

$ws.Range('D8').Value = '3.505.00'
$ws.Range('E8').Value = '  -0.88%  '

$ws.Range('E9').Value = '  -0.01%  '

$ws.Range('D10').Value = '''0.195'
$ws.Range('E10').Value = '  -2.66%  '

$ws.Range('E11').Value = '  +6.58%  '

$ws.Range('E12').Value = '  -0.09%  '

$ws.Range('D13').Value = '''46.21'
$ws.Range('E13').Value = '  -3.20%  '

$ws.Range('D14').Value = '''0.0000277'
$ws.Range('E14').Value = '  -1.67%  '

$ws.Range('D15').Value = '4.077.06'
$ws.Range('E15').Value = '  -0.90%  '

$ws.Range('D16').Value = '''8.37'
$ws.Range('E16').Value = '  -0.96%  '

$ws.Range('D17').Value = '''613.15'
$ws.Range('E17').Value = '  -2.69%  '

$ws.Range('D18').Value = '3.507.66'
$ws.Range('E18').Value = '  -0.45%  '

$ws.Range('D19').Value = '70.182.63'
$ws.Range('E19').Value = '  -0.47%  '

$ws.Range('E20').Value = '  +0.77%  '

$ws.Range('D21').Value = '''17.56'
$ws.Range('E21').Value = '  +0.50%  '

$ws.Range('D22').Value = '''0.878'
$ws.Range('E22').Value = '  -1.46%  '

$ws.Range('D23').Value = '''9.11'
$ws.Range('E23').Value = '  -10.08%  '

$ws.Range('D24').Value = '''98.78'
$ws.Range('E24').Value = '  +1.89%  '

$ws.Range('D25').Value = '''15.64'
$ws.Range('E25').Value = '  -1.68%  '

$ws.Range('D26').Value = '''3.72'
$ws.Range('E26').Value = '  -3.99%  '

$ws.Range('E27').Value = '  -0.02%  '

$ws.Range('D28').Value = '''2.57'
$ws.Range('E28').Value = '  -2.35%  '

$ws.Range('D29').Value = '''34.03'
$ws.Range('E29').Value = '  +1.76%  '

$ws.Range('E30').Value = '  -2.84%  '

$ws.Range('D31').Value = '''8.05'
$ws.Range('E31').Value = '  -4.96%  '

$ws.Range('D32').Value = '''2.97'
$ws.Range('E32').Value = '  -4.64%  '

$ws.Range('E33').Value = '  -5.11%  '

$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').Value = '''6.83'
$ws.Range('E34').Value = '  -3.23%  '

$ws.Range('B35').Value = 'Bittensor'
$ws.Range('C35').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D35').Value = '''630.08'
$ws.Range('E35').Value = '  +11.01%  '

$ws.Range('D36').Value = '''0.0996'

$ws.Range('D37').Value = '''10.75'
$ws.Range('E37').Value = '  -0.69%  '

$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '''0.0482'
$ws.Range('E38').Value = '  +5.88%  '

$ws.Range('B39').Value = 'dogwifhat'
$ws.Range('C39').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D39').Value = '''3.49'
$ws.Range('E39').Value = '  -3.03%  '

$ws.Range('D40').Value = '''56.76'
$ws.Range('E40').Value = '  -1.14%  '

$ws.Range('E41').Value = '  +0.19%  '

$ws.Range('D42').Value = '''0.144'
$ws.Range('E42').Value = '  +0.61%  '

$ws.Range('D43').Value = '3.367.95'
$ws.Range('E43').Value = '  +0.71%  '

$ws.Range('D44').Value = '0.0₃0733'
$ws.Range('E44').Value = '  +1.32%  '

$ws.Range('D45').Value = '''0.310'
$ws.Range('E45').Value = '  -6.07%  '

$ws.Range('D46').Value = '''2.92'
$ws.Range('E46').Value = '  -4.65%  '

$ws.Range('D47').Value = '''31.87'
$ws.Range('E47').Value = '  -4.19%  '

$ws.Range('D48').Value = '''2.55'
$ws.Range('E48').Value = '  -4.10%  '

$ws.Range('E49').Value = '  +0.07%  '

$ws.Range('D50').Value = '''133.09'
$ws.Range('E50').Value = '  -1.45%  '

$ws.Range('E51').Value = '  -0.03%  '
